# Generate Report for Handoff
# Updates status text from "In Translation" to "Ready for handoff" across
# the Overview, zh-cn, and de-de sheets, refreshes the related timestamp
# cells, and widens the status columns to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update status values ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Update timestamps (text values, not real dates, matching source data) ---
$overview.Range("G2").Value = "2016-08-31 18:46:09"
$dede.Range("H2").Value = "2016-08-31 18:46:09"
$zhcn.Range("H2").Value = "2016-08-31 18:45:57"

# --- Widen the status columns to fit the new, longer text ---
# (the new text "Ready for handoff" is wider than "In Translation",
# so Excel's autofit grows these columns)
$overview.Range("E:E").ColumnWidth = 16.3
$overview.Range("F:F").ColumnWidth = 16.3
$zhcn.Range("C:C").ColumnWidth = 16.3
$dede.Range("C:C").ColumnWidth = 16.3
